function Split-RunBoundary {
    param($d, $Start, $End)
    # Toggling a character property on/off over a span forces the engine
    # to materialize that span as its own <w:r> (formatting round-trips
    # back to identical, but the run boundary sticks) instead of merging
    # it back into the neighboring run.
    $rng = $d.Range($Start, $End)
    $rng.Bold = 1
    $rng.Bold = 0
}

function Insert-SplitRun {
    param($d, $InsertAt, $Text)
    $r = $d.Range($InsertAt, $InsertAt)
    $r.InsertAfter($Text)
    $txtLen = $Text.Length
    $endAt = $InsertAt + $txtLen
    Split-RunBoundary $d $InsertAt $endAt
}

$d = $word.ActiveDocument

# 1) "Title" -> "Title" + new run " 2"
$p = $d.Paragraphs.Item(1)
$insertAt = $p.Range.End - 1
Insert-SplitRun $d $insertAt " 2"

# 2) "Header 1" -> "Header 1" + new run "-2"
$p = $d.Paragraphs.Item(2)
$insertAt = $p.Range.End - 1
Insert-SplitRun $d $insertAt "-2"

# 3) "Header 2" -> "Header 2-2" (text changed in place, same run)
$d.Content.Find.Execute("Header 2", $true, $false, $false, $false, $false, $true, 1, $false, "Header 2-2", 2) | Out-Null

# 4) "Header 3" -> "Header 3" + new run "-2"
$p = $d.Paragraphs.Item(6)
$insertAt = $p.Range.End - 1
Insert-SplitRun $d $insertAt "-2"

# 5) "The end." -> "The end" + new run " 2" + new run "."
$p = $d.Paragraphs.Item(9)
$dotPos = $p.Range.End - 2   # position of "." (just before the paragraph mark)
$r = $d.Range($dotPos, $dotPos)
$r.InsertAfter(" 2")
$afterInsert = $dotPos + 2
Split-RunBoundary $d $dotPos $afterInsert
$afterDot = $afterInsert + 1
Split-RunBoundary $d $afterInsert $afterDot

Write-Output "Done"
